# Weekly update for Hortaliza, Vega Central Mapocho de Santiago - Rabanito
# Shifts existing rows 29-141 (new data inserted, existing rows pushed down) and
# appends trailing rows through row 148.
# Columns: Row, D(date serial), J(Volumen), K(Precio minimo), L(Precio maximo),
#          M(Precio promedio ponderado), O(Origen), P(Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(29,44453,7900,3500,4000,3747,"Provincia de Chacabuco",37),
    @(30,44162,20000,2500,3000,2725,"Provincia de Chacabuco",27),
    @(31,44413,8800,2500,3000,2750,"Provincia de Chacabuco",28),
    @(32,44202,8000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(33,44323,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(34,44169,21000,2500,3000,2762,"Provincia de Chacabuco",28),
    @(35,44216,12000,2500,3000,2708,"Provincia de Chacabuco",27),
    @(36,44405,5200,2500,3000,2750,"Provincia de Chacabuco",28),
    @(37,44260,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(38,44281,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(39,44308,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(40,44271,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(41,44363,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(42,44229,13000,2500,3000,2769,"Provincia de Chacabuco",28),
    @(43,44189,13000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(44,44181,14000,2500,3000,2714,"Provincia de Chacabuco",27),
    @(45,44231,11000,2500,3000,2727,"Provincia de Chacabuco",27),
    @(46,44355,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(47,44246,7000,3000,3000,3000,"Región Metropolitana",30),
    @(48,44203,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(49,44230,11000,2500,3000,2727,"Provincia de Chacabuco",27),
    @(50,44418,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(51,44258,10000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(52,44407,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(53,44236,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(54,44159,17000,2500,3000,2735,"Provincia de Chacabuco",27),
    @(55,44245,16000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(56,44330,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(57,44328,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(58,44187,19000,2500,3000,2763,"Provincia de Chacabuco",28),
    @(59,44343,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(60,44293,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(61,44344,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(62,44195,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(63,44320,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(64,44209,13000,2500,3000,2731,"Provincia de Chacabuco",27),
    @(65,44278,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(66,44406,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(67,44385,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(68,44250,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(69,44221,5000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(70,44334,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(71,44299,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(72,44265,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(73,44186,4000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(74,44252,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(75,44371,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(76,44419,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(77,44370,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(78,44447,7000,3500,4000,3750,"Provincia de Chacabuco",38),
    @(79,44267,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(80,44292,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(81,44259,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(82,44300,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(83,44392,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(84,44217,11000,2500,3000,2727,"Provincia de Chacabuco",27),
    @(85,44295,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(86,44362,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(87,44309,8000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(88,44384,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(89,44168,21000,2500,3000,2738,"Provincia de Chacabuco",27),
    @(90,44322,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(91,44272,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(92,44196,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(93,44365,6000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(94,44358,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(95,44218,13000,2500,3000,2692,"Provincia de Chacabuco",27),
    @(96,44433,7900,3500,4000,3750,"Provincia de Chacabuco",38),
    @(97,44397,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(98,44274,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(99,44321,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(100,44291,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(101,44335,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(102,44426,6100,3000,3500,3250,"Provincia de Chacabuco",32),
    @(103,44421,9700,2500,3000,2750,"Provincia de Chacabuco",28),
    @(104,44434,7900,3500,4000,3750,"Provincia de Chacabuco",38),
    @(105,44215,16000,2500,3000,2812,"Provincia de Chacabuco",28),
    @(106,44194,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(107,44222,15000,2500,3000,2767,"Provincia de Chacabuco",28),
    @(108,44398,7000,2500,3000,2750,"Provincia de Chacabuco",28),
    @(109,44420,9700,2500,3000,2750,"Provincia de Chacabuco",28),
    @(110,44264,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(111,44316,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(112,44210,17000,2500,3000,2765,"Provincia de Chacabuco",28),
    @(113,44253,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(114,44427,7000,3000,3500,3250,"Provincia de Chacabuco",32),
    @(115,44341,9700,3000,3000,3000,"Provincia de Chacabuco",30),
    @(116,44414,7900,2500,3000,2750,"Provincia de Chacabuco",28),
    @(117,44301,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(118,44176,14000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(119,44237,10000,2500,3000,2700,"Provincia de Chacabuco",27),
    @(120,44432,7000,3500,4000,3750,"Provincia de Chacabuco",38),
    @(121,44351,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(122,44369,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(123,44273,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(124,44342,9000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(125,44294,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(126,44379,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(127,44302,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(128,44315,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(129,44391,4300,2500,3000,2750,"Región Metropolitana",28),
    @(130,44446,7900,3000,4000,3494,"Provincia de Chacabuco",35),
    @(131,44411,6900,2500,3000,2750,"Provincia de Chacabuco",28),
    @(132,44313,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(133,44329,12000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(134,44449,7900,3500,4000,3753,"Provincia de Chacabuco",38),
    @(135,44161,14000,2500,3000,2679,"Provincia de Chacabuco",27),
    @(136,44251,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(137,44435,30700,3500,4000,3750,"Provincia de Chacabuco",38),
    @(138,44175,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(139,44376,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(140,44279,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(141,44412,7900,2500,3000,2782,"Provincia de Chacabuco",28),
    @(142,44223,13000,2500,3000,2808,"Provincia de Chacabuco",28),
    @(143,44314,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(144,44448,7900,3500,4000,3747,"Provincia de Chacabuco",37),
    @(145,44167,15000,2500,3000,2733,"Provincia de Chacabuco",27),
    @(146,44238,7000,3000,3000,3000,"Provincia de Chacabuco",30),
    @(147,44399,8800,2500,3000,2750,"Provincia de Chacabuco",28),
    @(148,44400,7000,2500,3000,2750,"Provincia de Chacabuco",28)
)

# Static values shared by every data row in this block (unchanged by the edit)
$colA = 9
$colB = "Vega Central Mapocho de Santiago"
$colC = "Metropolitana"
$colE = 13
$colF = 300000001
$colG = "Rabanito"
$colH = "Sin especificar"
$colI = "Primera"
$colN = "`$/cien unidades (volumen en unidades)"
$colQ = 100
$colR = "Hortaliza"

$lastExistingRow = 141
$newLastRow = 148

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $dateSerial = $entry[1]
    $vol = $entry[2]
    $pmin = $entry[3]
    $pmax = $entry[4]
    $pprom = $entry[5]
    $origen = $entry[6]
    $pkg = $entry[7]

    if ($r -gt $lastExistingRow) {
        # Brand new row beyond the previous table extent: populate every column.
        $ws.Cells.Item($r, 1).Value = $colA
        $ws.Cells.Item($r, 2).Value = $colB
        $ws.Cells.Item($r, 3).Value = $colC
        $ws.Cells.Item($r, 4).Value = $dateSerial
        $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Cells.Item($r, 5).Value = $colE
        $ws.Cells.Item($r, 6).Value = $colF
        $ws.Cells.Item($r, 7).Value = $colG
        $ws.Cells.Item($r, 8).Value = $colH
        $ws.Cells.Item($r, 9).Value = $colI
        $ws.Cells.Item($r, 10).Value = $vol
        $ws.Cells.Item($r, 11).Value = $pmin
        $ws.Cells.Item($r, 12).Value = $pmax
        $ws.Cells.Item($r, 13).Value = $pprom
        $ws.Cells.Item($r, 14).Value = $colN
        $ws.Cells.Item($r, 15).Value = $origen
        $ws.Cells.Item($r, 16).Value = $pkg
        $ws.Cells.Item($r, 17).Value = $colQ
        $ws.Cells.Item($r, 18).Value = $colR
    } else {
        # Existing row: only the fields affected by the weekly shift change.
        $ws.Cells.Item($r, 4).Value = $dateSerial
        $ws.Cells.Item($r, 10).Value = $vol
        $ws.Cells.Item($r, 11).Value = $pmin
        $ws.Cells.Item($r, 12).Value = $pmax
        $ws.Cells.Item($r, 13).Value = $pprom
        $ws.Cells.Item($r, 15).Value = $origen
        $ws.Cells.Item($r, 16).Value = $pkg
    }
}

Write-Output "Done updating rows 29-148"
